$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'259.13"
$ws.Range("E2").Value = "'0.52%"

$ws.Range("D3").Value = "'27.00"
$ws.Range("E3").Value = "'-0.59%"

$ws.Range("D4").Value = "'4.681"
$ws.Range("E4").Value = "'0.36%"

$ws.Range("D5").Value = "'0.06020"
$ws.Range("E5").Value = "'2.24%"

$ws.Range("D6").Value = "'6.677"
$ws.Range("E6").Value = "'0.53%"

$ws.Range("D7").Value = "'0.8592"
$ws.Range("E7").Value = "'0.17%"

$ws.Range("D8").Value = "'0.9254"
$ws.Range("E8").Value = "'-4.33%"

$ws.Range("D9").Value = "'0.1400"
$ws.Range("E9").Value = "'-0.54%"

$ws.Range("D10").Value = "'0.04928"
$ws.Range("E10").Value = "'23.66%"

$ws.Range("D11").Value = "'0.07027"
$ws.Range("E11").Value = "'-0.95%"

$ws.Range("D12").Value = "'0.03123"
$ws.Range("E12").Value = "'-1.73%"

$ws.Range("D13").Value = "'0.09117"
$ws.Range("E13").Value = "'-0.61%"

$ws.Range("D14").Value = "'0.001533"
$ws.Range("E14").Value = "'-1.12%"

$ws.Range("D15").Value = "'0.0006054"
$ws.Range("E15").Value = "'-94.26%"

$ws.Range("D16").Value = "'0.006002"
$ws.Range("E16").Value = "'-3.51%"

$ws.Range("E17").Value = "'-1.53%"

$ws.Range("D19").Value = "'2.166"
$ws.Range("E19").Value = "'-1.72%"

$ws.Range("E20").Value = "'0.44%"

$ws.Range("D21").Value = "'0.1299"
$ws.Range("E21").Value = "'0.48%"

$ws.Range("D22").Value = "'4.126"
$ws.Range("E22").Value = "'7.00%"

$ws.Range("D23").Value = "'0.04227"
$ws.Range("E23").Value = "'-0.01%"

$ws.Range("D24").Value = "'0.001217"
$ws.Range("E24").Value = "'-0.27%"

$ws.Range("D25").Value = "'0.004038"

$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'0.04%"

$ws.Range("D40").Value = "'0.03843"
$ws.Range("E40").Value = "'0.39%"

$ws.Range("D41").Value = "'0.1114"
$ws.Range("E41").Value = "'0.99%"

$ws.Range("D42").Value = "'0.003977"
$ws.Range("E42").Value = "'-36.66%"

$ws.Range("D43").Value = "'0.01504"
$ws.Range("E43").Value = "'31.45%"

$ws.Range("D44").Value = "'0.002200"
$ws.Range("E44").Value = "'0.03%"

$ws.Range("D45").Value = "'0.00005108"
$ws.Range("E45").Value = "'-6.43%"

$ws.Range("E46").Value = "'0.03%"

$ws.Range("D47").Value = "'0.05459"
$ws.Range("E47").Value = "'-9.01%"

$ws.Range("D48").Value = "'0.1353"
$ws.Range("E48").Value = "'4.48%"

$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'0.03%"

$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'0.03%"
